# Auto commit - 11101027
# Applies the swapped / corrected values for the row pairs that were
# re-sequenced in the source report (work-order numbers 2025110902/903,
# 2025110907/908, 2025110670/671, 2025110666/667, 2025110900/901).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

function Set-RowValues {
    param($Row, $B, $O, $S, $AA, $AB)

    $ws.Cells.Item($Row, 2).Value = $B          # B - 工作序號
    $ws.Cells.Item($Row, 15).Value = $O         # O - 工作類型

    if ($S -eq $null) {
        $ws.Cells.Item($Row, 19).Value = $null  # S - 保養
    } else {
        $ws.Cells.Item($Row, 19).Value = $S
    }

    if ($AA -eq $null) {
        $ws.Cells.Item($Row, 27).Value = $null  # AA - 合約事宜
    } else {
        $ws.Cells.Item($Row, 27).Value = $AA
    }

    if ($AB -eq $null) {
        $ws.Cells.Item($Row, 28).Value = $null  # AB - 其它
    } else {
        $ws.Cells.Item($Row, 28).Value = $AB
    }
}

# Row 48 / 49 pair
Set-RowValues 48 2025110903 "抄表" $null "O" 1
Set-RowValues 49 2025110902 "服務" "O" $null $null

# Row 58 / 59 pair (S only toggles; AA/AB already equal on both rows)
Set-RowValues 58 2025110908 "抄表" $null "O" 1
Set-RowValues 59 2025110907 "服務" "O" "O" $null

# Row 75 / 76 pair
Set-RowValues 75 2025110671 "抄表" $null "O" 1
Set-RowValues 76 2025110670 "服務" "O" $null $null

# Row 80 / 81 pair
Set-RowValues 80 2025110666 "服務" "O" $null $null
Set-RowValues 81 2025110667 "抄表" $null "O" $null

# Row 87 / 88 pair
Set-RowValues 87 2025110900 "服務" "O" $null $null
Set-RowValues 88 2025110901 "抄表" $null "O" $null
